# Auto-generated edit script applying cell-value changes described by the diff.
# Each entry: sheet name, cell reference, new value ($null means clear the cell).
$wb = $excel.ActiveWorkbook

$changes = @{
  "ALC" = @(
    @{ Cell = "H19"; Value = 2222.6667 }
    @{ Cell = "I19"; Value = 656.3333 }
    @{ Cell = "J19"; Value = 2744.7778 }
    @{ Cell = "K19"; Value = 656.3333 }
    @{ Cell = "L19"; Value = 2744.7778 }
    @{ Cell = "M19"; Value = -481.3333 }
    @{ Cell = "N19"; Value = -3094.7778 }
    @{ Cell = "H39"; Value = 1162.375 }
    @{ Cell = "J39"; Value = 1314.0714 }
    @{ Cell = "L39"; Value = 3942.2142 }
    @{ Cell = "N39"; Value = -4534.2142 }
    @{ Cell = "H40"; Value = 2236.889 }
    @{ Cell = "I40"; Value = 2225.5 }
    @{ Cell = "J40"; Value = 2259.6667 }
    @{ Cell = "K40"; Value = 2225.5 }
    @{ Cell = "L40"; Value = 2259.6667 }
    @{ Cell = "M40"; Value = -2050.5 }
    @{ Cell = "N40"; Value = -2609.6667 }
    @{ Cell = "H43"; Value = 1510.5834 }
    @{ Cell = "I43"; Value = 1573.4 }
    @{ Cell = "J43"; Value = 1465.7142 }
    @{ Cell = "K43"; Value = 1573.4 }
    @{ Cell = "L43"; Value = 1465.7142 }
    @{ Cell = "M43"; Value = -1504.4 }
    @{ Cell = "N43"; Value = -1603.7142 }
    @{ Cell = "H55"; Value = 189 }
    @{ Cell = "I55"; Value = 185 }
    @{ Cell = "K55"; Value = 185 }
    @{ Cell = "M55"; Value = 29 }
    @{ Cell = "H135"; Value = 756.4706 }
    @{ Cell = "I135"; Value = 479.22223 }
    @{ Cell = "J135"; Value = 1068.375 }
    @{ Cell = "K135"; Value = 4313.00007 }
    @{ Cell = "L135"; Value = 9615.375 }
    @{ Cell = "M135"; Value = -1778.00007 }
    @{ Cell = "N135"; Value = -14685.375 }
    @{ Cell = "H137"; Value = 2377.25 }
    @{ Cell = "J137"; Value = 2377.25 }
    @{ Cell = "L137"; Value = 7131.75 }
    @{ Cell = "N137"; Value = -12231.75 }
    @{ Cell = "H138"; Value = 2262.5144 }
    @{ Cell = "I138"; Value = 2262.2163 }
    @{ Cell = "J138"; Value = 2262.8484 }
    @{ Cell = "K138"; Value = 6786.6489 }
    @{ Cell = "L138"; Value = 6788.5452 }
    @{ Cell = "M138"; Value = -1646.6489 }
    @{ Cell = "N138"; Value = -17068.5452 }
  )
  "ARM" = @(
    @{ Cell = "H2"; Value = 278440.6 }
    @{ Cell = "I2"; Value = 347912.25 }
    @{ Cell = "J2"; Value = 554 }
    @{ Cell = "K2"; Value = 347912.25 }
    @{ Cell = "L2"; Value = 554 }
    @{ Cell = "M2"; Value = -347799.25 }
    @{ Cell = "N2"; Value = -780 }
    @{ Cell = "H32"; Value = 7287.51 }
    @{ Cell = "I32"; Value = 5043.268 }
    @{ Cell = "J32"; Value = 16488.9 }
    @{ Cell = "K32"; Value = 5043.268 }
    @{ Cell = "L32"; Value = 16488.9 }
    @{ Cell = "M32"; Value = -4756.268 }
    @{ Cell = "N32"; Value = -17062.9 }
    @{ Cell = "H38"; Value = 18000 }
    @{ Cell = "I38"; Value = 18000 }
    @{ Cell = "K38"; Value = 18000 }
    @{ Cell = "M38"; Value = -17533 }
    @{ Cell = "H74"; Value = 859.35 }
    @{ Cell = "I74"; Value = 565.02856 }
    @{ Cell = "K74"; Value = 565.02856 }
    @{ Cell = "M74"; Value = 308.97144 }
    @{ Cell = "H77"; Value = 859.35 }
    @{ Cell = "I77"; Value = 565.02856 }
    @{ Cell = "K77"; Value = 2825.1428 }
    @{ Cell = "M77"; Value = 1542.8572 }
    @{ Cell = "H116"; Value = 278440.6 }
    @{ Cell = "I116"; Value = 347912.25 }
    @{ Cell = "J116"; Value = 554 }
    @{ Cell = "K116"; Value = 347912.25 }
    @{ Cell = "L116"; Value = 554 }
    @{ Cell = "M116"; Value = -345618.25 }
    @{ Cell = "N116"; Value = -5142 }
    @{ Cell = "H122"; Value = 2999.875 }
    @{ Cell = "I122"; Value = 1250 }
    @{ Cell = "K122"; Value = 3750 }
    @{ Cell = "M122"; Value = -1300 }
  )
  "BSM" = @(
    @{ Cell = "H3"; Value = 278440.6 }
    @{ Cell = "I3"; Value = 347912.25 }
    @{ Cell = "J3"; Value = 554 }
    @{ Cell = "K3"; Value = 347912.25 }
    @{ Cell = "L3"; Value = 554 }
    @{ Cell = "M3"; Value = -347798.25 }
    @{ Cell = "N3"; Value = -782 }
    @{ Cell = "H38"; Value = 20000 }
    @{ Cell = "J38"; Value = 20000 }
    @{ Cell = "L38"; Value = 20000 }
    @{ Cell = "N38"; Value = -20832 }
    @{ Cell = "H129"; Value = 44999.7 }
    @{ Cell = "J129"; Value = 46666.332 }
    @{ Cell = "L129"; Value = 46666.332 }
    @{ Cell = "N129"; Value = -56666.332 }
  )
  "CRP" = @(
    @{ Cell = "H16"; Value = 764.8182 }
    @{ Cell = "I16"; Value = 686.2857 }
    @{ Cell = "J16"; Value = 902.25 }
    @{ Cell = "K16"; Value = 686.2857 }
    @{ Cell = "L16"; Value = 902.25 }
    @{ Cell = "M16"; Value = -399.2857 }
    @{ Cell = "N16"; Value = -1476.25 }
    @{ Cell = "H31"; Value = 3350.6667 }
    @{ Cell = "I31"; Value = 4071.5 }
    @{ Cell = "J31"; Value = 2990.25 }
    @{ Cell = "K31"; Value = 4071.5 }
    @{ Cell = "L31"; Value = 2990.25 }
    @{ Cell = "M31"; Value = -3776.5 }
    @{ Cell = "N31"; Value = -3580.25 }
    @{ Cell = "H34"; Value = 3350.6667 }
    @{ Cell = "I34"; Value = 4071.5 }
    @{ Cell = "J34"; Value = 2990.25 }
    @{ Cell = "K34"; Value = 4071.5 }
    @{ Cell = "L34"; Value = 2990.25 }
    @{ Cell = "M34"; Value = -3869.5 }
    @{ Cell = "N34"; Value = -3394.25 }
    @{ Cell = "H35"; Value = 2541.6667 }
    @{ Cell = "I35"; Value = 1312.5 }
    @{ Cell = "J35"; Value = 5000 }
    @{ Cell = "K35"; Value = 1312.5 }
    @{ Cell = "L35"; Value = 5000 }
    @{ Cell = "M35"; Value = -1018.5 }
    @{ Cell = "N35"; Value = -5588 }
    @{ Cell = "H38"; Value = 1900 }
    @{ Cell = "I38"; Value = 1900 }
    @{ Cell = "J38"; Value = 0 }
    @{ Cell = "K38"; Value = 1900 }
    @{ Cell = "L38"; Value = 0 }
    @{ Cell = "M38"; Value = -1523 }
    @{ Cell = "N38"; Value = $null }
    @{ Cell = "H46"; Value = 1900 }
    @{ Cell = "I46"; Value = 1900 }
    @{ Cell = "J46"; Value = 0 }
    @{ Cell = "K46"; Value = 1900 }
    @{ Cell = "L46"; Value = 0 }
    @{ Cell = "M46"; Value = -1689 }
    @{ Cell = "N46"; Value = $null }
    @{ Cell = "H113"; Value = 764.8182 }
    @{ Cell = "I113"; Value = 686.2857 }
    @{ Cell = "J113"; Value = 902.25 }
    @{ Cell = "K113"; Value = 686.2857 }
    @{ Cell = "L113"; Value = 902.25 }
    @{ Cell = "M113"; Value = 1483.7143 }
    @{ Cell = "N113"; Value = -5242.25 }
    @{ Cell = "H134"; Value = 1841.4651 }
    @{ Cell = "I134"; Value = 1701.0541 }
    @{ Cell = "K134"; Value = 5103.1623 }
    @{ Cell = "M134"; Value = -2568.1623 }
  )
  "CUL" = @(
    @{ Cell = "H2"; Value = 210.25 }
    @{ Cell = "I2"; Value = 147.22223 }
    @{ Cell = "J2"; Value = 399.33334 }
    @{ Cell = "K2"; Value = 883.33338 }
    @{ Cell = "L2"; Value = 2396.00004 }
    @{ Cell = "M2"; Value = -770.33338 }
    @{ Cell = "N2"; Value = -2622.00004 }
    @{ Cell = "H11"; Value = 861.4 }
    @{ Cell = "I11"; Value = 770 }
    @{ Cell = "K11"; Value = 2310 }
    @{ Cell = "M11"; Value = -2170 }
    @{ Cell = "H33"; Value = 666850.0600000001 }
    @{ Cell = "I33"; Value = 120.42857 }
    @{ Cell = "J33"; Value = 1250238.5 }
    @{ Cell = "K33"; Value = 722.57142 }
    @{ Cell = "L33"; Value = 7501431 }
    @{ Cell = "M33"; Value = -439.57142 }
    @{ Cell = "N33"; Value = -7501997 }
    @{ Cell = "H131"; Value = 25676.322 }
    @{ Cell = "J131"; Value = 28696.32 }
    @{ Cell = "L131"; Value = 86088.95999999999 }
    @{ Cell = "N131"; Value = -96168.95999999999 }
    @{ Cell = "H132"; Value = 1372.5834 }
    @{ Cell = "I132"; Value = 1020.6 }
    @{ Cell = "K132"; Value = 9185.4 }
    @{ Cell = "M132"; Value = -6655.4 }
  )
  "GSM" = @(
    @{ Cell = "H113"; Value = 930.1111 }
    @{ Cell = "I113"; Value = 660.7 }
    @{ Cell = "K113"; Value = 660.7 }
    @{ Cell = "M113"; Value = 1509.3 }
    @{ Cell = "H132"; Value = 1242887.4 }
    @{ Cell = "I132"; Value = 1833241 }
    @{ Cell = "J132"; Value = 3144.8 }
    @{ Cell = "K132"; Value = 5499723 }
    @{ Cell = "L132"; Value = 9434.400000000001 }
    @{ Cell = "M132"; Value = -5497193 }
    @{ Cell = "N132"; Value = -14494.4 }
  )
  "LTW" = @(
    @{ Cell = "H16"; Value = 5234.385 }
    @{ Cell = "I16"; Value = 5234.385 }
    @{ Cell = "K16"; Value = 5234.385 }
    @{ Cell = "M16"; Value = -5064.385 }
    @{ Cell = "H38"; Value = 10000 }
    @{ Cell = "J38"; Value = 10000 }
    @{ Cell = "L38"; Value = 10000 }
    @{ Cell = "N38"; Value = -10820 }
    @{ Cell = "H61"; Value = 2713.1143 }
    @{ Cell = "I61"; Value = 2431.6296 }
    @{ Cell = "K61"; Value = 2431.6296 }
    @{ Cell = "M61"; Value = -2229.6296 }
    @{ Cell = "H63"; Value = 46985 }
    @{ Cell = "J63"; Value = 46985 }
    @{ Cell = "L63"; Value = 46985 }
    @{ Cell = "N63"; Value = -48483 }
    @{ Cell = "H66"; Value = 46985 }
    @{ Cell = "J66"; Value = 46985 }
    @{ Cell = "L66"; Value = 140955 }
    @{ Cell = "N66"; Value = -148443 }
    @{ Cell = "H74"; Value = 10000 }
    @{ Cell = "I74"; Value = 10000 }
    @{ Cell = "J74"; Value = 0 }
    @{ Cell = "K74"; Value = 10000 }
    @{ Cell = "L74"; Value = 0 }
    @{ Cell = "M74"; Value = -9002 }
    @{ Cell = "N74"; Value = $null }
    @{ Cell = "H77"; Value = 10000 }
    @{ Cell = "I77"; Value = 10000 }
    @{ Cell = "J77"; Value = 0 }
    @{ Cell = "K77"; Value = 30000 }
    @{ Cell = "L77"; Value = 0 }
    @{ Cell = "M77"; Value = -25008 }
    @{ Cell = "N77"; Value = $null }
    @{ Cell = "H82"; Value = 2858.6667 }
    @{ Cell = "I82"; Value = 1956 }
    @{ Cell = "J82"; Value = 3987 }
    @{ Cell = "K82"; Value = 1956 }
    @{ Cell = "L82"; Value = 3987 }
    @{ Cell = "M82"; Value = -1595 }
    @{ Cell = "N82"; Value = -4709 }
    @{ Cell = "H85"; Value = 2858.6667 }
    @{ Cell = "I85"; Value = 1956 }
    @{ Cell = "J85"; Value = 3987 }
    @{ Cell = "K85"; Value = 1956 }
    @{ Cell = "L85"; Value = 3987 }
    @{ Cell = "M85"; Value = -708 }
    @{ Cell = "N85"; Value = -6483 }
    @{ Cell = "H113"; Value = 2713.1143 }
    @{ Cell = "I113"; Value = 2431.6296 }
    @{ Cell = "K113"; Value = 2431.6296 }
    @{ Cell = "M113"; Value = -261.6296000000002 }
    @{ Cell = "H122"; Value = 26999.572 }
    @{ Cell = "I122"; Value = 27499.5 }
    @{ Cell = "K122"; Value = 82498.5 }
    @{ Cell = "M122"; Value = -80048.5 }
  )
  "WVR" = @(
    @{ Cell = "H62"; Value = 6111 }
    @{ Cell = "I62"; Value = 3333 }
    @{ Cell = "K62"; Value = 3333 }
    @{ Cell = "M62"; Value = -2709 }
    @{ Cell = "H65"; Value = 6111 }
    @{ Cell = "I65"; Value = 3333 }
    @{ Cell = "K65"; Value = 16665 }
    @{ Cell = "M65"; Value = -13545 }
    @{ Cell = "H126"; Value = 6734.476 }
    @{ Cell = "I126"; Value = 10194.667 }
    @{ Cell = "K126"; Value = 30584.001 }
    @{ Cell = "M126"; Value = -28114.001 }
    @{ Cell = "H132"; Value = 1841.7307 }
    @{ Cell = "I132"; Value = 1270.409 }
    @{ Cell = "J132"; Value = 4984 }
    @{ Cell = "K132"; Value = 3811.227 }
    @{ Cell = "L132"; Value = 14952 }
    @{ Cell = "M132"; Value = -1281.227 }
    @{ Cell = "N132"; Value = -20012 }
  )
}

foreach ($sheetName in $changes.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)
  foreach ($chg in $changes[$sheetName]) {
    $rng = $ws.Range($chg.Cell)
    if ($null -eq $chg.Value) {
      $rng.ClearContents()
    } else {
      $rng.Value = $chg.Value
    }
  }
}

Write-Host "Applied $(($changes.Values | ForEach-Object { $_.Count } | Measure-Object -Sum).Sum) cell changes across $($changes.Keys.Count) sheets."